$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.855.22'
$ws.Range("E2").Value = '  -4.09%  '

$ws.Range("D3").Value = '3.504.34'
$ws.Range("E3").Value = '  -4.81%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.75'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.76%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.10'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -3.20%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.622'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.20%  '

$ws.Range("D8").Value = '3.495.08'
$ws.Range("E8").Value = '  -4.89%  '

$ws.Range("E9").Value = '  +0.07%  '

$ws.Range("E10").Value = '  -7.31%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.54'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +4.38%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.603'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.79%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '47.20'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -5.60%  '

$ws.Range("E14").Value = '  -3.85%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '675.05'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.26%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.89'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.50%  '

$ws.Range("D17").Value = '4.063.48'
$ws.Range("E17").Value = '  -4.68%  '

$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.501.91'
$ws.Range("E18").Value = '  -4.69%  '

$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = '68.819.71'
$ws.Range("E19").Value = '  -4.29%  '

$ws.Range("E20").Value = '  -1.76%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.56'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -4.11%  '

$ws.Range("E22").Value = '  -4.65%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.905'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -4.31%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '16.30'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -8.72%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '98.24'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -5.33%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.87'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -4.29%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.81'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.93%  '

$ws.Range("E28").Value = '  +0.09%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.66'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -6.57%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.43'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -8.15%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.98'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -7.10%  '

$ws.Range("E32").Value = '  -5.09%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.21'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -7.50%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.35'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -1.47%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.35'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -6.43%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '571.91'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.64%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.61'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -15.04%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '10.94'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -3.73%  '

$ws.Range("E39").Value = '  -3.56%  '

$ws.Range("E40").Value = '  -5.78%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.02%  '

$ws.Range("E42").Value = '  -4.99%  '

$ws.Range("E43").Value = '  -4.93%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.337'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -3.24%  '

$ws.Range("D45").Value = '3.420.94'
$ws.Range("E45").Value = '  -8.64%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '33.40'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -6.49%  '

$ws.Range("D47").Value = '0.0₃0703'
$ws.Range("E47").Value = '  -9.28%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.89'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.98%  '

$ws.Range("E49").Value = '  -7.69%  '

$ws.Range("E50").Value = '  -0.82%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '133.58'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.47%  '
